$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new AlexNet row values (row 8, columns B-E) with percentage formatted numbers
$ws.Range("B8").Value = 0.6118
$ws.Range("C8").Value = 0.8292
$ws.Range("D8").Value = 0.6057
$ws.Range("E8").Value = 0.8236

$ws.Range("B8:E8").NumberFormat = "0.00%"
$ws.Range("B8:E8").HorizontalAlignment = -4108

# Update selection to G12
$ws.Range("G12").Select()
